# MAJ automatique BRVM - refresh des classements "Recommandations" / "Top_YTD"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet 1: "Recommandations" ---
# Sector/stock ranking reshuffled (rows re-sorted by updated "Variation Totale")
# plus refreshed Jours en Hausse/Baisse, Variation Totale, Derniere Variation,
# Recommandation and Strategie for every affected row.
$ws1.Cells.Item(2,4).Value = 2442.71
$ws1.Cells.Item(2,5).Value = 102.52
$ws1.Cells.Item(3,4).Value = 2102.31
$ws1.Cells.Item(3,5).Value = 685.18
$ws1.Cells.Item(4,4).Value = 1490.97
$ws1.Cells.Item(4,5).Value = 496.73
$ws1.Cells.Item(5,4).Value = 1084.35
$ws1.Cells.Item(5,5).Value = 369.61
$ws1.Cells.Item(6,4).Value = 995.34
$ws1.Cells.Item(6,5).Value = 331.47
$ws1.Cells.Item(7,1).Value = 'BRVM - INDUSTRIE  (**)'
$ws1.Cells.Item(7,3).Value = 3
$ws1.Cells.Item(7,4).Value = 767.38
$ws1.Cells.Item(7,5).Value = 255.06
$ws1.Cells.Item(8,1).Value = 'BRVM-PRINCIPAL  (**)'
$ws1.Cells.Item(8,4).Value = 648.64
$ws1.Cells.Item(8,5).Value = 215.96
$ws1.Cells.Item(9,1).Value = 'BRVM - CONSOMMATION DE BASE  (**)'
$ws1.Cells.Item(9,4).Value = 639.99
$ws1.Cells.Item(9,5).Value = 212.87
$ws1.Cells.Item(10,1).Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Cells.Item(10,4).Value = 527.61
$ws1.Cells.Item(10,5).Value = 174.53
$ws1.Cells.Item(11,1).Value = 'BRVM - FINANCES'
$ws1.Cells.Item(11,4).Value = 439.58
$ws1.Cells.Item(11,5).Value = 147.59
$ws1.Cells.Item(12,1).Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Cells.Item(12,4).Value = 432.01
$ws1.Cells.Item(12,5).Value = 145.05
$ws1.Cells.Item(13,1).Value = 'BRVM-PRESTIGE'
$ws1.Cells.Item(13,4).Value = 423.91
$ws1.Cells.Item(13,5).Value = 142.16
$ws1.Cells.Item(14,1).Value = 'BRVM - INDUSTRIELS'
$ws1.Cells.Item(14,4).Value = 397.49
$ws1.Cells.Item(14,5).Value = 132.17
$ws1.Cells.Item(15,1).Value = 'BRVM - ENERGIE'
$ws1.Cells.Item(15,4).Value = 331.19
$ws1.Cells.Item(15,5).Value = 111.03
$ws1.Cells.Item(16,1).Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Cells.Item(16,4).Value = 278.04
$ws1.Cells.Item(16,5).Value = 91.74
$ws1.Cells.Item(17,1).Value = 'NEI-CEDA CI (NEIC)'
$ws1.Cells.Item(17,2).Value = 3
$ws1.Cells.Item(17,3).Value = 0
$ws1.Cells.Item(17,4).Value = 20.06
$ws1.Cells.Item(17,5).Value = 7.21
$ws1.Cells.Item(17,6).Value = '🟢 Achat'
$ws1.Cells.Item(17,7).Value = '✅ Renforcer'
$ws1.Cells.Item(18,2).Value = 2
$ws1.Cells.Item(18,3).Value = 1
$ws1.Cells.Item(18,4).Value = 7.32
$ws1.Cells.Item(18,5).Value = -7.19
$ws1.Cells.Item(18,6).Value = '🟡 Observer'
$ws1.Cells.Item(18,7).Value = '👀 À surveiller'
$ws1.Cells.Item(19,1).Value = 'SICOR CI (SICC)'
$ws1.Cells.Item(19,2).Value = 1
$ws1.Cells.Item(19,4).Value = 6.71
$ws1.Cells.Item(19,5).Value = 6.71
$ws1.Cells.Item(20,1).Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws1.Cells.Item(20,2).Value = 1
$ws1.Cells.Item(20,4).Value = 4.75
$ws1.Cells.Item(20,5).Value = 4.75
$ws1.Cells.Item(21,1).Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Cells.Item(21,4).Value = 4.13
$ws1.Cells.Item(21,5).Value = 4.13
$ws1.Cells.Item(22,1).Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Cells.Item(22,4).Value = 3.7
$ws1.Cells.Item(22,5).Value = 3.7
$ws1.Cells.Item(23,1).Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Cells.Item(23,4).Value = 3.42
$ws1.Cells.Item(23,5).Value = 3.42
$ws1.Cells.Item(24,1).Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = 2.9
$ws1.Cells.Item(24,5).Value = -3.65
$ws1.Cells.Item(24,7).Value = '👀 À surveiller'
$ws1.Cells.Item(25,1).Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Cells.Item(25,4).Value = 0.81
$ws1.Cells.Item(25,5).Value = -3.43
$ws1.Cells.Item(26,1).Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Cells.Item(26,2).Value = 1
$ws1.Cells.Item(26,3).Value = 1
$ws1.Cells.Item(26,4).Value = 0.2
$ws1.Cells.Item(26,5).Value = 4.55
$ws1.Cells.Item(26,7).Value = '👀 À surveiller'
$ws1.Cells.Item(27,1).Value = 'TOTAL'
$ws1.Cells.Item(27,3).Value = 2
$ws1.Cells.Item(27,4).Value = 0
$ws1.Cells.Item(27,5).Value = 0
$ws1.Cells.Item(28,1).Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,4).Value = -0.04
$ws1.Cells.Item(28,5).Value = 3.96
$ws1.Cells.Item(28,7).Value = '👀 À surveiller'
$ws1.Cells.Item(29,1).Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,4).Value = -1.16
$ws1.Cells.Item(29,5).Value = -6.02
$ws1.Cells.Item(29,7).Value = '👀 À surveiller'
$ws1.Cells.Item(30,1).Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws1.Cells.Item(30,2).Value = 0
$ws1.Cells.Item(30,4).Value = -1.5
$ws1.Cells.Item(30,5).Value = -1.5
$ws1.Cells.Item(30,7).Value = '➖ Neutre'
$ws1.Cells.Item(31,1).Value = 'NESTLE CI (NTLC)'
$ws1.Cells.Item(31,3).Value = 2
$ws1.Cells.Item(31,4).Value = -3.74
$ws1.Cells.Item(31,5).Value = -2.78
$ws1.Cells.Item(32,1).Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Cells.Item(32,2).Value = 0
$ws1.Cells.Item(32,3).Value = 1
$ws1.Cells.Item(32,4).Value = -4.92
$ws1.Cells.Item(32,5).Value = -4.92
$ws1.Cells.Item(32,7).Value = '➖ Neutre'
$ws1.Cells.Item(33,1).Value = 'SUCRIVOIRE (SCRC)'
$ws1.Cells.Item(33,4).Value = -6.35
$ws1.Cells.Item(33,5).Value = -5.45
$ws1.Cells.Item(34,1).Value = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$ws1.Cells.Item(34,3).Value = 3
$ws1.Cells.Item(34,4).Value = -16.33
$ws1.Cells.Item(34,5).Value = -6.25
$ws1.Cells.Item(34,6).Value = '🔴 Vente'
$ws1.Cells.Item(34,7).Value = '⚠️ Risque de décrochage'

# Rows 35:36 no longer exist in the updated table -> clear them
$ws1.Range("A35:G36").ClearContents()

# --- Sheet 2: "Top_YTD" ---
$ws2.Cells.Item(2,2).Value = 448888.34
$ws2.Cells.Item(3,2).Value = 51228.93
$ws2.Cells.Item(4,2).Value = 21176.43
$ws2.Cells.Item(5,2).Value = 9723.65
$ws2.Cells.Item(6,2).Value = 7949.82
$ws2.Cells.Item(7,2).Value = 4403.89
$ws2.Cells.Item(8,2).Value = 3061.84
$ws2.Cells.Item(9,2).Value = 2976.12
$ws2.Cells.Item(10,2).Value = 1999.44
$ws2.Cells.Item(11,2).Value = 1398.25
